$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: current_phase 1 -> 2
$ws.Range("D20").Value = 2

# Row 21: current_phase 1 -> 2
$ws.Range("D21").Value = 2

# Row 21: last_action_date set
$ws.Range("E21").Value = "2026-02-12T13:19:54.271941+00:00"

# Row 21: reactions_count 0 -> 1
$ws.Range("H21").Value = 1

# Row 21: replies_count 0 -> 1
$ws.Range("I21").Value = 1

# Row 21: reacted_message_ids [] -> [19]
$ws.Range("L21").Value = "[19]"

# Row 21: replied_message_ids [] -> [27]
$ws.Range("M21").Value = "[27]"
